$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the width (COM character units) of column G before inserting,
# so the newly inserted column H can inherit a matching width.
$gColumnWidth = $ws.Columns.Item(7).ColumnWidth()

# Insert a new column before H (shifts H:Z to I:AA).
$ws.Range("H1").EntireColumn.Insert()

# Apply the width that column G had, to the freshly inserted column H.
$ws.Columns.Item(8).ColumnWidth = $gColumnWidth

# New header cell for the inserted column.
$ws.Range("H1").Value = "Work Location"

# Restore the view selection similarly to what the authored change shows.
$ws.Range("Z15").Select()
